$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "C7"  = -12.965
    "B9"  = 5.317
    "C12" = -11.171
    "B18" = 5.282999999999999
    "B20" = 6.952
    "C26" = -13.131
    "B27" = 5.526
    "C27" = -13.664
    "C29" = -12.335
    "C37" = -13.426
    "C38" = -13.313
    "C51" = -11.312
    "C55" = -13.752
    "B69" = 5.667
    "C69" = -11.17
    "C70" = -11.645
    "B76" = 6.723999999999999
    "B82" = 5.241
    "C83" = -13.425
    "C102" = -13.419
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
